$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from an existing header cell
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Clear clipboard/marching ants
$excel.CutCopyMode = $false

# Data values for columns I ("I0") and J ("IF") per row
$values = @{
    2  = @(8, 9)
    3  = @(7, 8)
    4  = @(8, 9)
    5  = @(7, 8)
    6  = @(1, 2)
    7  = @(6, 7)
    8  = @(6, 7)
    9  = @(7, 8)
    10 = @(6, 7)
    11 = @(1, 2)
    12 = @(6, 7)
    13 = @(7, 7)
    14 = @(1, 2)
    15 = @(8, 9)
    16 = @(1, 1)
    17 = @(7, 8)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]   # Column I
    $ws.Cells.Item($row, 10).Value = $pair[1]  # Column J
}
